# "tolte dispense modifica DT dati"
# - Shared strings: remove the old "lotto" entry and add a new "Lotto"
#   entry (capitalised) at the end of the table; this is what the header
#   cell B1 on Foglio1 points at.
# - Update the selection on Foglio1 from the whole-sheet default
#   (A1:D1048576) to a normal single active cell at B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")
$ws.Activate()

# Re-label the "lotto" header as "Lotto" (new shared-string entry).
$ws.Range("B1").Value = "Lotto"

# Collapse the selection down to B2 (was the full-column A1:D1048576).
$ws.Range("B2").Select()
